$wb = $excel.ActiveWorkbook

# --- "Branch LPF" sheet: cutoff freq 4800 -> 2400, sample rate 28800 -> 14400, gain 2 -> 1 ---
$wsBranch = $wb.Worksheets.Item("Branch LPF")
$wsBranch.Range("B2").Value = 2400
$wsBranch.Range("B3").Value = 14400
$wsBranch.Range("B18").Value = 1
[void]$wsBranch.Range("B23").Select()

# --- "LoopFilter LPF" sheet: cutoff freq 100 -> 50, sample rate 28800 -> 14400, gain 1 -> 8 ---
$wsLoop = $wb.Worksheets.Item("LoopFilter LPF")
$wsLoop.Range("B2").Value = 50
$wsLoop.Range("B3").Value = 14400
$wsLoop.Range("B16").Value = 8
[void]$wsLoop.Range("B17").Select()
